$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Declined")

$ws.Range("A2").Value = "2025-08-06T01:53:22.718Z"
$ws.Range("B2").Value = "ton nguyen"

$ws.Range("A3").Value = "2025-08-06T01:54:07.017Z"
$ws.Range("B3").Value = "ton nguyen"

$ws.Range("A4").Value = "2025-08-06T01:54:08.635Z"
$ws.Range("B4").Value = "ton nguyen"
